# "added I as a word"
#
# The speech originally opened directly with the "Graduation Speech"
# title paragraph. The edit inserts a brand-new first paragraph whose
# only content is the single word "I", pushing the title (and
# everything after it) down by one paragraph.
#
# Grab the very first paragraph in the document (currently "Graduation
# Speech") and splice a new paragraph in front of it containing just
# "I", followed by a paragraph mark (carriage return) so it becomes its
# own, separate paragraph rather than merging into the title line.
$d = $word.ActiveDocument
$firstParagraph = $d.Paragraphs(1).Range
$firstParagraph.InsertBefore("I" + [char]13)
